$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A165").Value = "2023-12-10 16:17:21"
$ws.Range("B165").Value = 0.0004

$ws.Range("A166").Value = "2023-12-10 16:17:26"
$ws.Range("B166").Value = 0.0006000000000000001
